# Auto-generated Excel COM-interop script applying numeric corrections
# to columns H-N across multiple rows on several worksheets, per the
# scheduled-runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 864.5217
$ws.Range("I18").Value = 623.0476
$ws.Range("J18").Value = 3400
$ws.Range("K18").Value = 623.0476
$ws.Range("L18").Value = 3400
$ws.Range("M18").Value = -339.0476
$ws.Range("N18").Value = -3968
$ws.Range("H132").Value = 4812229.5
$ws.Range("I132").Value = 5439483.5
$ws.Range("J132").Value = 3285
$ws.Range("K132").Value = 16318450.5
$ws.Range("L132").Value = 9855
$ws.Range("M132").Value = -16315920.5
$ws.Range("N132").Value = -14915
$ws.Range("H137").Value = 1463.279
$ws.Range("I137").Value = 1039.9445
$ws.Range("J137").Value = 3640.4285
$ws.Range("K137").Value = 3119.8335
$ws.Range("L137").Value = 10921.2855
$ws.Range("M137").Value = -569.8335000000002
$ws.Range("N137").Value = -16021.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29804.92
$ws.Range("I32").Value = 6042.5293
$ws.Range("K32").Value = 6042.5293
$ws.Range("M32").Value = -5755.5293
$ws.Range("H45").Value = 42557.48
$ws.Range("I45").Value = 73640.28999999999
$ws.Range("J45").Value = 2997.5454
$ws.Range("K45").Value = 73640.28999999999
$ws.Range("L45").Value = 2997.5454
$ws.Range("M45").Value = -73263.28999999999
$ws.Range("N45").Value = -3751.5454
$ws.Range("H61").Value = 1508.44
$ws.Range("I61").Value = 1405.2858
$ws.Range("J61").Value = 2050
$ws.Range("K61").Value = 1405.2858
$ws.Range("L61").Value = 2050
$ws.Range("M61").Value = -1193.2858
$ws.Range("N61").Value = -2474
$ws.Range("H74").Value = 1459.8667
$ws.Range("I74").Value = 747.7143
$ws.Range("J74").Value = 3121.5557
$ws.Range("K74").Value = 747.7143
$ws.Range("L74").Value = 3121.5557
$ws.Range("M74").Value = 126.2857
$ws.Range("N74").Value = -4869.5557
$ws.Range("H77").Value = 1459.8667
$ws.Range("I77").Value = 747.7143
$ws.Range("J77").Value = 3121.5557
$ws.Range("K77").Value = 3738.5715
$ws.Range("L77").Value = 15607.7785
$ws.Range("M77").Value = 629.4285
$ws.Range("N77").Value = -24343.7785
$ws.Range("H88").Value = 2755.5557
$ws.Range("I88").Value = 2575
$ws.Range("J88").Value = 2900
$ws.Range("K88").Value = 2575
$ws.Range("L88").Value = 2900
$ws.Range("M88").Value = -2169
$ws.Range("N88").Value = -3712
$ws.Range("H91").Value = 2755.5557
$ws.Range("I91").Value = 2575
$ws.Range("J91").Value = 2900
$ws.Range("K91").Value = 2575
$ws.Range("L91").Value = 2900
$ws.Range("M91").Value = -1171
$ws.Range("N91").Value = -5708
$ws.Range("H132").Value = 2072.9355
$ws.Range("I132").Value = 1286.8422
$ws.Range("J132").Value = 3317.5833
$ws.Range("K132").Value = 3860.5266
$ws.Range("L132").Value = 9952.749899999999
$ws.Range("M132").Value = -1330.5266
$ws.Range("N132").Value = -15012.7499
$ws.Range("H136").Value = 1508.44
$ws.Range("I136").Value = 1405.2858
$ws.Range("J136").Value = 2050
$ws.Range("K136").Value = 4215.857400000001
$ws.Range("L136").Value = 6150
$ws.Range("M136").Value = -1665.857400000001
$ws.Range("N136").Value = -11250

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 41125.465
$ws.Range("I86").Value = 53818.383
$ws.Range("J86").Value = 3046.7144
$ws.Range("K86").Value = 53818.383
$ws.Range("L86").Value = 3046.7144
$ws.Range("M86").Value = -52695.383
$ws.Range("N86").Value = -5292.7144
$ws.Range("H89").Value = 41125.465
$ws.Range("I89").Value = 53818.383
$ws.Range("J89").Value = 3046.7144
$ws.Range("K89").Value = 269091.915
$ws.Range("L89").Value = 15233.572
$ws.Range("M89").Value = -263475.915
$ws.Range("N89").Value = -26465.572
$ws.Range("H134").Value = 3209.8
$ws.Range("I134").Value = 3338.2307
$ws.Range("J134").Value = 2375
$ws.Range("K134").Value = 10014.6921
$ws.Range("L134").Value = 7125
$ws.Range("M134").Value = -7479.6921
$ws.Range("N134").Value = -12195

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 65889.27
$ws.Range("I31").Value = 59597.234
$ws.Range("J31").Value = 70539.914
$ws.Range("K31").Value = 59597.234
$ws.Range("L31").Value = 70539.914
$ws.Range("M31").Value = -59302.234
$ws.Range("N31").Value = -71129.914
$ws.Range("H34").Value = 65889.27
$ws.Range("I34").Value = 59597.234
$ws.Range("J34").Value = 70539.914
$ws.Range("K34").Value = 59597.234
$ws.Range("L34").Value = 70539.914
$ws.Range("M34").Value = -59395.234
$ws.Range("N34").Value = -70943.914
$ws.Range("H41").Value = 10726
$ws.Range("I41").Value = 4275
$ws.Range("J41").Value = 15026.667
$ws.Range("K41").Value = 4275
$ws.Range("L41").Value = 15026.667
$ws.Range("M41").Value = -3847
$ws.Range("N41").Value = -15882.667
$ws.Range("H50").Value = 11710
$ws.Range("J50").Value = 11710
$ws.Range("L50").Value = 11710
$ws.Range("N50").Value = -12960
$ws.Range("H51").Value = 7874.9287
$ws.Range("I51").Value = 8090
$ws.Range("J51").Value = 7858.385
$ws.Range("K51").Value = 8090
$ws.Range("L51").Value = 7858.385
$ws.Range("M51").Value = -7354
$ws.Range("N51").Value = -9330.385
$ws.Range("H58").Value = 5781.8613
$ws.Range("I58").Value = 1465.1364
$ws.Range("J58").Value = 12565.286
$ws.Range("K58").Value = 1465.1364
$ws.Range("L58").Value = 12565.286
$ws.Range("M58").Value = -1262.1364
$ws.Range("N58").Value = -12971.286
$ws.Range("H59").Value = 22717.273
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 22717.273
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 22717.273
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -25007.273
$ws.Range("H60").Value = 10252.074
$ws.Range("I60").Value = 8052
$ws.Range("J60").Value = 10428.08
$ws.Range("K60").Value = 8052
$ws.Range("L60").Value = 10428.08
$ws.Range("M60").Value = -7541
$ws.Range("N60").Value = -11450.08
$ws.Range("H61").Value = 7874.9287
$ws.Range("I61").Value = 8090
$ws.Range("J61").Value = 7858.385
$ws.Range("K61").Value = 8090
$ws.Range("L61").Value = 7858.385
$ws.Range("M61").Value = -7742
$ws.Range("N61").Value = -8554.385
$ws.Range("H68").Value = 15774.292
$ws.Range("J68").Value = 15774.292
$ws.Range("L68").Value = 15774.292
$ws.Range("N68").Value = -17272.292
$ws.Range("H71").Value = 15774.292
$ws.Range("J71").Value = 15774.292
$ws.Range("L71").Value = 47322.876
$ws.Range("N71").Value = -54810.876
$ws.Range("H74").Value = 24739.334
$ws.Range("J74").Value = 24739.334
$ws.Range("L74").Value = 24739.334
$ws.Range("N74").Value = -26487.334
$ws.Range("H77").Value = 24739.334
$ws.Range("J77").Value = 24739.334
$ws.Range("L77").Value = 74218.00199999999
$ws.Range("N77").Value = -82954.00199999999
$ws.Range("H96").Value = 22874.666
$ws.Range("J96").Value = 22874.666
$ws.Range("L96").Value = 22874.666
$ws.Range("N96").Value = -28366.666
$ws.Range("H132").Value = 2979.2183
$ws.Range("I132").Value = 2926.0557
$ws.Range("J132").Value = 3079.9473
$ws.Range("K132").Value = 8778.167099999999
$ws.Range("L132").Value = 9239.841899999999
$ws.Range("M132").Value = -6248.167099999999
$ws.Range("N132").Value = -14299.8419
$ws.Range("H134").Value = 1128.1063
$ws.Range("I134").Value = 1124.2413
$ws.Range("J134").Value = 1134.3334
$ws.Range("K134").Value = 3372.7239
$ws.Range("L134").Value = 3403.0002
$ws.Range("M134").Value = -837.7239
$ws.Range("N134").Value = -8473.0002
$ws.Range("H136").Value = 5781.8613
$ws.Range("I136").Value = 1465.1364
$ws.Range("J136").Value = 12565.286
$ws.Range("K136").Value = 4395.4092
$ws.Range("L136").Value = 37695.858
$ws.Range("M136").Value = -1845.4092
$ws.Range("N136").Value = -42795.858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 86.04000000000001
$ws.Range("I2").Value = 138.11111
$ws.Range("J2").Value = 56.75
$ws.Range("K2").Value = 828.66666
$ws.Range("L2").Value = 340.5
$ws.Range("M2").Value = -715.66666
$ws.Range("N2").Value = -566.5
$ws.Range("H5").Value = 1316.2894
$ws.Range("I5").Value = 619.76
$ws.Range("J5").Value = 2655.7693
$ws.Range("K5").Value = 1859.28
$ws.Range("L5").Value = 7967.3079
$ws.Range("M5").Value = -1747.28
$ws.Range("N5").Value = -8191.3079
$ws.Range("H34").Value = 437.9
$ws.Range("J34").Value = 761.8
$ws.Range("L34").Value = 2285.4
$ws.Range("N34").Value = -2453.4
$ws.Range("H39").Value = 2120
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 2525
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 7575
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -8163
$ws.Range("H55").Value = 8601.036
$ws.Range("I55").Value = 20378
$ws.Range("J55").Value = 6040.826
$ws.Range("K55").Value = 61134
$ws.Range("L55").Value = 18122.478
$ws.Range("M55").Value = -60957
$ws.Range("N55").Value = -18476.478
$ws.Range("H122").Value = 6822
$ws.Range("I122").Value = 473
$ws.Range("J122").Value = 34334.332
$ws.Range("K122").Value = 4257
$ws.Range("L122").Value = 309008.988
$ws.Range("M122").Value = -1807
$ws.Range("N122").Value = -313908.988
$ws.Range("H131").Value = 857.39
$ws.Range("J131").Value = 857.39
$ws.Range("L131").Value = 2572.17
$ws.Range("N131").Value = -12652.17
$ws.Range("H135").Value = 1316.2894
$ws.Range("I135").Value = 619.76
$ws.Range("J135").Value = 2655.7693
$ws.Range("K135").Value = 5577.84
$ws.Range("L135").Value = 23901.9237
$ws.Range("M135").Value = -3042.84
$ws.Range("N135").Value = -28971.9237

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3756.8572
$ws.Range("I132").Value = 2433.3333
$ws.Range("J132").Value = 4749.5
$ws.Range("K132").Value = 7299.999899999999
$ws.Range("L132").Value = 14248.5
$ws.Range("M132").Value = -4769.999899999999
$ws.Range("N132").Value = -19308.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2838.158
$ws.Range("J132").Value = 2181
$ws.Range("L132").Value = 6543
$ws.Range("N132").Value = -11603
$ws.Range("H136").Value = 1447.4
$ws.Range("I136").Value = 1386.125
$ws.Range("J136").Value = 1692.5
$ws.Range("K136").Value = 4158.375
$ws.Range("L136").Value = 5077.5
$ws.Range("M136").Value = -1608.375
$ws.Range("N136").Value = -10177.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 37495
$ws.Range("J127").Value = 37495
$ws.Range("L127").Value = 37495
$ws.Range("N127").Value = -47415
$ws.Range("H132").Value = 2634.7
$ws.Range("I132").Value = 2114.4102
$ws.Range("J132").Value = 4479.364
$ws.Range("K132").Value = 6343.230599999999
$ws.Range("L132").Value = 13438.092
$ws.Range("M132").Value = -3813.230599999999
$ws.Range("N132").Value = -18498.092
$ws.Range("H136").Value = 1110.1818
$ws.Range("I136").Value = 926.5
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 2779.5
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = -229.5
$ws.Range("N136").Value = -9900
